$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "julia.tyndan7@gmail.com"
$ws.Range("B3").Value = "simiyu.wdan32@gmail.com"

$ws.Range("C9").Select()
